# Add a case when serviceInterface has annotations but serviceImpl doesn't have.
#
# Current layout (rows 1..11):
#   Row 1        : headers
#   Rows 2-3     : "В контроллере" cases
#   Rows 4-7     : "В сервисе" cases
#     Row 7      : "Сервис как интерфейс без аннотаций валидации + Impl c @Validated"
#   Row 8        : blank
#   Row 9        : "Вывод" / "В контроллере можно не ставить ..."
#   Row 10       : "В сервисах от аннотации `@Validated` на объектах нет смысла"
#   Row 11       : "Если сервис объявлен как интерфейс, то надо иметь валидационные
#                   аннотации и в интерфейсе и в классе имплементации"
#
# New layout adds a new table row (new row 8) describing the case where the
# service interface carries the validation annotations while the Impl class
# (annotated with @Validated) has none on its method parameters, pushes the
# "Вывод" block down by one row, and splits the old last bullet into two
# separate bullets (one kept in place, one new one appended at the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new table row by inserting a new row 8 (shifts the
#    "Вывод" block that used to start at row 9 down to row 10, etc.)
$ws.Rows("8:8").Insert()

# 2) Fill in the new table row (A8:E8) for the new case.
$ws.Range("A8").Value = "В сервисе"
$ws.Range("B8").Value = "Сервис как интерфейс с аннотациями валидации + Impl c @Validated, но без аннотаций на параметрах метода"
$ws.Range("C8").Value = "+ ConstraintViolationException"
$ws.Range("D8").Value = "+ ConstraintViolationException"
$ws.Range("E8").Value = "- не работает валидация"

# Style to match the other "В сервисе"-group first rows: green fill on A/B
# (same look as A5/B5, A7/B7) plus word-wrap on B8 because the text got
# longer, and a taller row to show the wrapped text; quote-prefixed,
# vertically centered cells on C8:E8 to match the existing quotePrefix style
# used throughout column C/D/E.
$ws.Range("A8").Style = "40% - Accent1"
$ws.Rows("8:8").RowHeight = 30

$ws.Range("A8:B8").Interior.Color = $ws.Range("A5:B5").Interior.Color
$ws.Range("B8").WrapText = $true
$ws.Range("C8:E8").VerticalAlignment = -4108

# 3) Extend the table ("Таблица1") and its AutoFilter to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E8"))

# 4) The old closing bullet (row 11 pre-insert, now row 12) explained that the
#    interface AND the impl class both needed validation annotations. Replace
#    it with the narrower statement (only the interface needs them) and add a
#    brand-new bullet after it clarifying the impl-class case we just added.
$ws.Range("B12").Value = "Если сервис объявлен как интерфейс, то надо иметь валидационные аннотации в интерфейсе (в классе имплементации не обязательно)"
$ws.Range("B13").Value = "Если сервис объявлен как интерфейс, то надо иметь валидационные аннотации в интерфейсе (в классе имплементации не обязательно)"
